# Laboratorio 8 - Entrega final
# Fill in student names and answer paragraphs to match the graded/final
# version of the lab report.

$d = $word.ActiveDocument

# --- Student 1 name / code ------------------------------------------------
$d.Content.Find.Execute(
    "Estudiante 1 Cod XXXX", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Mario Alberto Ricaurte Tobón; 201922994", 2) | Out-Null

# --- Student 2 name / code ------------------------------------------------
$d.Content.Find.Execute(
    "Estudiante 2 Cod XXXX", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Laura Valentina Jiménez Tobar; 201924116", 2) | Out-Null

# --- Answer to question 1 (the empty paragraph right after it) -----------
# "¿Qué relación encuentra entre el número de elementos en el árbol y la
#  altura del árbol?" -> paragraph 6 is the first (still empty) paragraph
#  that follows it.
$p1 = $d.Paragraphs.Item(6)
$p1.Range.InsertXML('<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="34F36316" w14:textId="23F61F9F" w:rsidR="0063268C" w:rsidRDefault="0063268C" w:rsidP="0063268C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:after="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:lang w:val="es-CO"/></w:rPr><w:t>Encontramos que la relación entre el número de elementos y la altura es del orden 29</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:vertAlign w:val="superscript"/><w:lang w:val="es-CO"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> puesto que es un árbol binario, no obstante, este resultado es menor que la cantidad de datos. Esto se debe a que el árbol está desbalanceado. Encontramos que puede estar desbalanceado porque la raíz no es un número central entre todas las fechas. Además, también está desbalanceado porque hay fechas en las que ocurren más crímenes que en otras, por ejemplo, en días de fiesta o durante los fines de semana. Esta correlación también desbalancea el árbol.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Answer to question 2 (the empty paragraph right after it) -----------
# "¿Si tuviera que responder esa misma consulta ... tablas de hash ...?"
# -> paragraph 9 is the first (still empty) paragraph that follows it.
$p2 = $d.Paragraphs.Item(9)
$p2.Range.InsertXML('<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="67B54B66" w14:textId="7ACD5B0C" w:rsidR="0063268C" w:rsidRDefault="0063268C" w:rsidP="0063268C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:lang w:val="es-CO"/></w:rPr><w:t>Tomaría más tiempo, porque en una tabla de Hash tendría que recorrer toda la tabla y tomando los datos que estén en la fecha específica. Esto significa que tendría que hacer una comparación por cada dato en la tabla, más aún si para solucionar colisiones se crea una lista dentro de cada llave, se harían más comparaciones. En el árbol es más rápido, puesto a que ya está ordenado por fechas, entonces se va descartando la mitad del árbol en cada comparación, teniendo una complejidad de ln(N).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Answer to question 3 (the empty paragraph right after it) -----------
# "¿Qué operación del TAD se utiliza para retornar una lista ...?"
# -> paragraph 12 is the first (still empty) paragraph that follows it.
$p3 = $d.Paragraphs.Item(12)
$p3.Range.InsertXML('<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1FAAE2DA" w14:textId="77777777" w:rsidR="0063268C" w:rsidRPr="0063268C" w:rsidRDefault="0063268C" w:rsidP="0063268C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">Se utiliza la operación </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:i/><w:iCs/><w:lang w:val="es-CO"/></w:rPr><w:t>values(map, keylo, keyhi)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:lang w:val="es-CO"/></w:rPr><w:t>. Esta operación recibe el mapa ordenado, es decir el árbol, un valor de la llave mínimo y un valor de la llave máximo, y retorna todos los hijos comprendidos entre esas dos llaves.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

Write-Output "done"
